# Updates TPM-derived NATMI statistics for Cxcl13-Ccr10 LR pairs sheet.
# Target cluster for some rows is corrected from "MuSCs" to "Resolving-Mac",
# and downstream specificity/weight metrics are recomputed with the new TPM data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    # Row 2
    $ws.Range("I2").Value = 0.5162107379131895
    $ws.Range("J2").Value = 0.5162107379131895
    $ws.Range("K2").Value = 3
    $ws.Range("L2").Value = 1
    $ws.Range("M2").Value = 1.058783666666667
    $ws.Range("N2").Value = 3.176351
    $ws.Range("O2").Value = 0.5014862149947701
    $ws.Range("P2").Value = 0.5014862149947702
    $ws.Range("Q2").Value = 4.284331402666222
    $ws.Range("R2").Value = 38.558982623996
    $ws.Range("S2").Value = 0.2588725690957427
    $ws.Range("T2").Value = 0.2588725690957427
    # Row 3
    $ws.Range("I3").Value = 0.5162107379131895
    $ws.Range("J3").Value = 0.5162107379131895
    $ws.Range("O3").Value = 0.475327031240749
    $ws.Range("P3").Value = 0.4753270312407492
    $ws.Range("S3").Value = 0.2453689175468727
    $ws.Range("T3").Value = 0.2453689175468728
    # Row 4
    $ws.Range("D4").Value = "Resolving-Mac"
    $ws.Range("I4").Value = 0.5162107379131895
    $ws.Range("J4").Value = 0.5162107379131895
    $ws.Range("M4").Value = 0.048954
    $ws.Range("N4").Value = 0.146862
    $ws.Range("O4").Value = 0.02318675376448066
    $ws.Range("P4").Value = 0.02318675376448067
    $ws.Range("Q4").Value = 0.198090663928
    $ws.Range("R4").Value = 1.782815975352
    $ws.Range("S4").Value = 0.01196925127057399
    $ws.Range("T4").Value = 0.01196925127057399
    # Row 5
    $ws.Range("G5").Value = 0.155986
    $ws.Range("H5").Value = 0.467958
    $ws.Range("I5").Value = 0.01989925565426652
    $ws.Range("J5").Value = 0.01989925565426652
    $ws.Range("K5").Value = 3
    $ws.Range("L5").Value = 1
    $ws.Range("M5").Value = 1.058783666666667
    $ws.Range("N5").Value = 3.176351
    $ws.Range("O5").Value = 0.5014862149947701
    $ws.Range("P5").Value = 0.5014862149947702
    $ws.Range("Q5").Value = 0.1651554290286667
    $ws.Range("R5").Value = 1.486398861258
    $ws.Range("S5").Value = 0.009979202399271393
    $ws.Range("T5").Value = 0.009979202399271395
    # Row 6
    $ws.Range("G6").Value = 0.155986
    $ws.Range("H6").Value = 0.467958
    $ws.Range("I6").Value = 0.01989925565426652
    $ws.Range("J6").Value = 0.01989925565426652
    $ws.Range("O6").Value = 0.475327031240749
    $ws.Range("P6").Value = 0.4753270312407492
    $ws.Range("Q6").Value = 0.156540374244
    $ws.Range("R6").Value = 1.408863368196
    $ws.Range("S6").Value = 0.009458654114043194
    $ws.Range("T6").Value = 0.009458654114043195
    # Row 7
    $ws.Range("D7").Value = "Resolving-Mac"
    $ws.Range("G7").Value = 0.155986
    $ws.Range("H7").Value = 0.467958
    $ws.Range("I7").Value = 0.01989925565426652
    $ws.Range("J7").Value = 0.01989925565426652
    $ws.Range("M7").Value = 0.048954
    $ws.Range("N7").Value = 0.146862
    $ws.Range("O7").Value = 0.02318675376448066
    $ws.Range("P7").Value = 0.02318675376448067
    $ws.Range("Q7").Value = 0.007636138643999999
    $ws.Range("R7").Value = 0.06872524779599999
    $ws.Range("S7").Value = 0.0004613991409519273
    $ws.Range("T7").Value = 0.0004613991409519274
    # Row 8
    $ws.Range("E8").Value = 3
    $ws.Range("F8").Value = 1
    $ws.Range("G8").Value = 3.636334333333334
    $ws.Range("H8").Value = 10.909003
    $ws.Range("I8").Value = 0.463890006432544
    $ws.Range("J8").Value = 0.463890006432544
    $ws.Range("K8").Value = 3
    $ws.Range("L8").Value = 1
    $ws.Range("M8").Value = 1.058783666666667
    $ws.Range("N8").Value = 3.176351
    $ws.Range("O8").Value = 0.5014862149947701
    $ws.Range("P8").Value = 0.5014862149947702
    $ws.Range("Q8").Value = 3.850091398672556
    $ws.Range("R8").Value = 34.650822588053
    $ws.Range("S8").Value = 0.232634443499756
    $ws.Range("T8").Value = 0.2326344434997561
    # Row 9
    $ws.Range("E9").Value = 3
    $ws.Range("F9").Value = 1
    $ws.Range("G9").Value = 3.636334333333334
    $ws.Range("H9").Value = 10.909003
    $ws.Range("I9").Value = 0.463890006432544
    $ws.Range("J9").Value = 0.463890006432544
    $ws.Range("O9").Value = 0.475327031240749
    $ws.Range("P9").Value = 0.4753270312407492
    $ws.Range("Q9").Value = 3.649257865554
    $ws.Range("R9").Value = 32.843320789986
    $ws.Range("S9").Value = 0.2204994595798331
    $ws.Range("T9").Value = 0.2204994595798332
    # Row 10
    $ws.Range("D10").Value = "Resolving-Mac"
    $ws.Range("E10").Value = 3
    $ws.Range("F10").Value = 1
    $ws.Range("G10").Value = 3.636334333333334
    $ws.Range("H10").Value = 10.909003
    $ws.Range("I10").Value = 0.463890006432544
    $ws.Range("J10").Value = 0.463890006432544
    $ws.Range("M10").Value = 0.048954
    $ws.Range("N10").Value = 0.146862
    $ws.Range("O10").Value = 0.02318675376448066
    $ws.Range("P10").Value = 0.02318675376448067
    $ws.Range("Q10").Value = 0.178013110954
    $ws.Range("R10").Value = 1.602117998586
    $ws.Range("S10").Value = 0.01075610335295475
    $ws.Range("T10").Value = 0.01075610335295475
